$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4350.9
$ws.Range("I19").Value = 2211.375
$ws.Range("J19").Value = 5777.25
$ws.Range("K19").Value = 2211.375
$ws.Range("L19").Value = 5777.25
$ws.Range("M19").Value = -2036.375
$ws.Range("N19").Value = -6127.25

$ws.Range("H76").Value = 4295.5264
$ws.Range("I76").Value = 3976
$ws.Range("J76").Value = 5999.6665
$ws.Range("K76").Value = 3976
$ws.Range("L76").Value = 5999.6665
$ws.Range("M76").Value = -3661
$ws.Range("N76").Value = -6629.6665

$ws.Range("H79").Value = 4295.5264
$ws.Range("I79").Value = 3976
$ws.Range("J79").Value = 5999.6665
$ws.Range("K79").Value = 3976
$ws.Range("L79").Value = 5999.6665
$ws.Range("M79").Value = -2884
$ws.Range("N79").Value = -8183.6665

$ws.Range("H113").Value = 3077.8
$ws.Range("I113").Value = 2334.6
$ws.Range("J113").Value = 3821
$ws.Range("K113").Value = 2334.6
$ws.Range("L113").Value = 3821
$ws.Range("M113").Value = 919.4000000000001
$ws.Range("N113").Value = -10329

$ws.Range("H138").Value = 17472.338
$ws.Range("I138").Value = 1558.0817
$ws.Range("K138").Value = 4674.2451
$ws.Range("M138").Value = 465.7548999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 734.2727
$ws.Range("I5").Value = 734.2727
$ws.Range("K5").Value = 734.2727
$ws.Range("M5").Value = -622.2727

$ws.Range("H61").Value = 8214.532999999999
$ws.Range("I61").Value = 1101.5834
$ws.Range("K61").Value = 1101.5834
$ws.Range("M61").Value = -889.5834

$ws.Range("H74").Value = 872944.4399999999
$ws.Range("I74").Value = 1001768.5
$ws.Range("K74").Value = 1001768.5
$ws.Range("M74").Value = -1000894.5

$ws.Range("H77").Value = 872944.4399999999
$ws.Range("I77").Value = 1001768.5
$ws.Range("K77").Value = 5008842.5
$ws.Range("M77").Value = -5004474.5

$ws.Range("H88").Value = 6424.5454
$ws.Range("I88").Value = 2331.6667
$ws.Range("J88").Value = 7959.375
$ws.Range("K88").Value = 2331.6667
$ws.Range("L88").Value = 7959.375
$ws.Range("M88").Value = -1925.6667
$ws.Range("N88").Value = -8771.375

$ws.Range("H91").Value = 6424.5454
$ws.Range("I91").Value = 2331.6667
$ws.Range("J91").Value = 7959.375
$ws.Range("K91").Value = 2331.6667
$ws.Range("L91").Value = 7959.375
$ws.Range("M91").Value = -927.6667000000002
$ws.Range("N91").Value = -10767.375

$ws.Range("H110").Value = 1047.9412
$ws.Range("I110").Value = 787.2143
$ws.Range("J110").Value = 2264.6667
$ws.Range("K110").Value = 787.2143
$ws.Range("L110").Value = 2264.6667
$ws.Range("M110").Value = 1257.7857
$ws.Range("N110").Value = -6354.6667

$ws.Range("H132").Value = 889.8461
$ws.Range("I132").Value = 860.9459000000001
$ws.Range("K132").Value = 2582.8377
$ws.Range("M132").Value = -52.83770000000004

$ws.Range("H136").Value = 8214.532999999999
$ws.Range("I136").Value = 1101.5834
$ws.Range("K136").Value = 3304.7502
$ws.Range("M136").Value = -754.7501999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 734.2727
$ws.Range("I4").Value = 734.2727
$ws.Range("K4").Value = 734.2727
$ws.Range("M4").Value = -619.2727

$ws.Range("H86").Value = 1920
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 1920
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -21232

$ws.Range("H134").Value = 2082.0444
$ws.Range("I134").Value = 1598.1515
$ws.Range("K134").Value = 4794.4545
$ws.Range("M134").Value = -2259.4545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 927.4666999999999
$ws.Range("J22").Value = 1433.3334
$ws.Range("L22").Value = 1433.3334
$ws.Range("N22").Value = -2133.3334

$ws.Range("H23").Value = 12000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 12000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H62").Value = 6281.857
$ws.Range("I62").Value = 6749.5
$ws.Range("J62").Value = 5658.3335
$ws.Range("K62").Value = 6749.5
$ws.Range("L62").Value = 5658.3335
$ws.Range("M62").Value = -6125.5
$ws.Range("N62").Value = -6906.3335

$ws.Range("H65").Value = 6281.857
$ws.Range("I65").Value = 6749.5
$ws.Range("J65").Value = 5658.3335
$ws.Range("K65").Value = 33747.5
$ws.Range("L65").Value = 28291.6675
$ws.Range("M65").Value = -30627.5
$ws.Range("N65").Value = -34531.6675

$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 102506.9
$ws.Range("I32").Value = 3338.3333
$ws.Range("J32").Value = 145007.72
$ws.Range("K32").Value = 10014.9999
$ws.Range("L32").Value = 435023.16
$ws.Range("M32").Value = -9731.999899999999
$ws.Range("N32").Value = -435589.16

$ws.Range("H139").Value = 10354.667
$ws.Range("I139").Value = 10354.667
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 31064.001
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -25924.001
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7701.722
$ws.Range("I70").Value = 7566.5713
$ws.Range("J70").Value = 8174.75
$ws.Range("K70").Value = 7566.5713
$ws.Range("L70").Value = 8174.75
$ws.Range("M70").Value = -7296.5713
$ws.Range("N70").Value = -8714.75

$ws.Range("H73").Value = 7701.722
$ws.Range("I73").Value = 7566.5713
$ws.Range("J73").Value = 8174.75
$ws.Range("K73").Value = 7566.5713
$ws.Range("L73").Value = 8174.75
$ws.Range("M73").Value = -6630.5713
$ws.Range("N73").Value = -10046.75

$ws.Range("H80").Value = 14122.5
$ws.Range("I80").Value = 1995
$ws.Range("J80").Value = 26250
$ws.Range("K80").Value = 1995
$ws.Range("L80").Value = 26250
$ws.Range("M80").Value = -997
$ws.Range("N80").Value = -28246

$ws.Range("H83").Value = 14122.5
$ws.Range("I83").Value = 1995
$ws.Range("J83").Value = 26250
$ws.Range("K83").Value = 9975
$ws.Range("L83").Value = 131250
$ws.Range("M83").Value = -4983
$ws.Range("N83").Value = -141234

$ws.Range("H94").Value = 9958.714
$ws.Range("J94").Value = 9958.714
$ws.Range("L94").Value = 9958.714
$ws.Range("N94").Value = -11310.714

$ws.Range("H140").Value = 79980
$ws.Range("J140").Value = 79980
$ws.Range("L140").Value = 79980
$ws.Range("N140").Value = -90340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2680.1875
$ws.Range("I46").Value = 1067.5
$ws.Range("K46").Value = 1067.5
$ws.Range("M46").Value = -879.5

$ws.Range("H82").Value = 1465.8667
$ws.Range("I82").Value = 1440.7142
$ws.Range("K82").Value = 1440.7142
$ws.Range("M82").Value = -1079.7142

$ws.Range("H85").Value = 1465.8667
$ws.Range("I85").Value = 1440.7142
$ws.Range("K85").Value = 1440.7142
$ws.Range("M85").Value = -192.7141999999999

$ws.Range("H93").Value = 2236.6365
$ws.Range("I93").Value = 2404
$ws.Range("J93").Value = 1943.75
$ws.Range("K93").Value = 2404
$ws.Range("L93").Value = 1943.75
$ws.Range("M93").Value = -1156
$ws.Range("N93").Value = -4439.75

$ws.Range("H132").Value = 3702.2222
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3702.2222
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11106.6666
$ws.Range("N132").Value = -16166.6666
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8782.723
$ws.Range("I81").Value = 10124.214
$ws.Range("J81").Value = 4087.5
$ws.Range("K81").Value = 20248.428
$ws.Range("L81").Value = 8175
$ws.Range("M81").Value = -19187.428
$ws.Range("N81").Value = -10297

$ws.Range("H84").Value = 8782.723
$ws.Range("I84").Value = 10124.214
$ws.Range("J84").Value = 4087.5
$ws.Range("K84").Value = 101242.14
$ws.Range("L84").Value = 40875
$ws.Range("M84").Value = -95938.14
$ws.Range("N84").Value = -51483

$ws.Range("H124").Value = 59999.5
$ws.Range("J124").Value = 59999.5
$ws.Range("L124").Value = 59999.5
$ws.Range("N124").Value = -69819.5
